# edit.ps1 - Apply "Doing Updates for Financials" quarterly refresh.
#
# Summary of the change (per the OOXML diff):
#  - Two new quarter columns are inserted before the old column D (so the
#    prior D:K data shifts right to F:M) on worksheet "HL".
#  - The two newly-inserted columns (D & E) are populated with the newest
#    two quarters of data for every metric row.
#  - A handful of cells within the shifted F:M range are restated with
#    corrected values (not pure shifts) - these are applied explicitly
#    after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HL")

# D/E column values (new 2 quarters) keyed by row -> (D,E)
$deData = @{
    7 = @(43465, 43373)
    8 = @(136500, 143600)
    9 = @(137800, 137100)
    10 = @(-1300, 6500)
    12 = @(9800, 14900)
    13 = @(0, 0)
    14 = @(2700, 12700)
    15 = @(0, 0)
    17 = @(161000, 174000)
    18 = @(-24500, -30400)
    20 = @(6500, 14700)
    21 = @(14100, 30200)
    22 = @(10900, 10100)
    23 = @(-28900, -25900)
    24 = @(-5200, -2700)
    25 = @(0, 0)
    26 = @(-23700, -23200)
    27 = @(-23800, -23300)
    28 = @(0, 0)
    29 = @("NA", "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-6500, -14700)
    33 = @(-23800, -23300)
    34 = @(0, 0)
    35 = @(-23800, -23300)
    38 = @(43465, 43373)
    41 = @(27400, 60900)
    42 = @(0, 0)
    43 = @(25800, 39900)
    44 = @(87500, 76100)
    45 = @(23400, 21500)
    46 = @(164200, 198300)
    47 = @(6600, 7200)
    48 = @(2520000, 2487400)
    49 = @(0, 0)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(13200, 17300)
    53 = @(0, 0)
    54 = @(2703900, 2710300)
    57 = @(77900, 65800)
    58 = @(5300, 6100)
    59 = @(53100, 60600)
    60 = @(136200, 132500)
    61 = @(540700, 542700)
    62 = @(336100, 313000)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(1013000, 988200)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(-248300, -223300)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(1690900, 1722000)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(-23800, -23300)
    83 = @(32100, 46000)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(19000, 28200)
    91 = @(-53600, -40000)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-52000, -169400)
    96 = @(-1200, -1200)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(800, -38000)
    101 = @(-1300, 300)
    102 = @(-33500, -178900)
}

# Corrections to shifted F:M cells that are not pure shifts (data restatements)
$corrections = @(
    ,@(8, "J", 276800)
    ,@(9, "H", 113800)
    ,@(9, "I", 97900)
    ,@(9, "J", 213600)
    ,@(10, "H", 46300)
    ,@(10, "I", 42900)
    ,@(10, "J", 63200)
    ,@(12, "J", 13700)
    ,@(14, "J", 9600)
    ,@(17, "H", 137200)
    ,@(17, "I", 120600)
    ,@(17, "J", 259900)
    ,@(18, "H", 22900)
    ,@(18, "I", 20200)
    ,@(18, "J", 16900)
    ,@(20, "H", -3800)
    ,@(20, "I", -15700)
    ,@(20, "J", -10100)
    ,@(21, "I", 38900)
    ,@(21, "J", 63700)
    ,@(22, "J", 19100)
    ,@(23, "H", 9600)
    ,@(23, "I", -4800)
    ,@(23, "J", -12300)
    ,@(24, "H", 68500)
    ,@(24, "I", -5100)
    ,@(24, "J", -12400)
    ,@(26, "H", -59000)
    ,@(26, "I", 300)
    ,@(26, "J", 100)
    ,@(27, "H", -59100)
    ,@(27, "I", 200)
    ,@(27, "J", -100)
    ,@(32, "H", 3800)
    ,@(32, "I", 15700)
    ,@(32, "J", 10100)
    ,@(33, "H", -29100)
    ,@(33, "I", 200)
    ,@(33, "J", -100)
    ,@(35, "H", -29100)
    ,@(35, "I", 200)
    ,@(35, "J", -100)
    ,@(44, "H", 55500)
    ,@(46, "H", 321200)
    ,@(48, "H", 1999300)
    ,@(54, "H", 2345200)
    ,@(62, "H", 263300)
    ,@(66, "H", 883900)
    ,@(72, "H", -218100)
    ,@(76, "H", 1461200)
    ,@(81, "H", -29100)
    ,@(81, "I", 200)
    ,@(81, "J", -100)
    ,@(83, "H", 35200)
    ,@(83, "I", 34300)
    ,@(83, "J", 56900)
    ,@(89, "J", 45800)
    ,@(91, "J", -46000)
    ,@(94, "J", -53600)
    ,@(96, "J", -2000)
    ,@(100, "J", 1100)
    ,@(101, "J", 1100)
    ,@(102, "J", -5700)
)

# ---------------------------------------------------------------------
# 1) Insert two new blank columns at D:E. Excel's EntireColumn.Insert()
#    shifts the existing D:K data (and its formatting) right to F:M.
# ---------------------------------------------------------------------
$ws.Range("D:E").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2) The freshly inserted D:E columns come in unformatted. Copy number
#    formatting/style down from column F (the original column D, now
#    shifted) row-by-row so the new columns visually match (date style
#    for the header row, numeric style for data rows, etc.).
# ---------------------------------------------------------------------
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Write the two newest quarters of data into columns D and E for
#    every populated row.
# ---------------------------------------------------------------------
foreach ($r in $deData.Keys) {
    $pair = $deData[$r]
    $ws.Cells.Item([int]$r, 4).Value = $pair[0]
    $ws.Cells.Item([int]$r, 5).Value = $pair[1]
}

# ---------------------------------------------------------------------
# 4) Apply restated values to specific cells in the shifted F:M range
#    that differ from a pure column shift (data corrections from the
#    source).
# ---------------------------------------------------------------------
$colMap = @{ "F" = 6; "G" = 7; "H" = 8; "I" = 9; "J" = 10; "K" = 11; "L" = 12; "M" = 13 }
foreach ($item in $corrections) {
    $r = $item[0]
    $colLetter = $item[1]
    $val = $item[2]
    $c = $colMap[$colLetter]
    $ws.Cells.Item([int]$r, $c).Value = $val
}
